$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.255.43'
$ws.Range("E2").Value = '  -0.50%  '
$ws.Range("D3").Value = '2.319.82'
$ws.Range("E3").Value = '  -1.77%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '533.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.58%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("E8").Value = '  -1.02%  '
$ws.Range("D9").Value = '2.344.83'
$ws.Range("E9").Value = '  -1.41%  '
$ws.Range("E10").Value = '  -1.51%  '
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("E12").Value = '  -2.81%  '
$ws.Range("E13").Value = '  -0.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.57'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.97%  '
$ws.Range("D15").Value = '2.740.00'
$ws.Range("E15").Value = '  -1.33%  '
$ws.Range("D16").Value = '57.271.02'
$ws.Range("E16").Value = '  -0.46%  '
$ws.Range("E17").Value = '  -2.44%  '
$ws.Range("D18").Value = '2.345.12'
$ws.Range("E18").Value = '  -0.71%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '340.13'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.45'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.30%  '
$ws.Range("E21").Value = '  +2.20%  '
$ws.Range("E22").Value = '  -2.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.79'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.36%  '
$ws.Range("E25").Value = '  +6.76%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.992'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.33'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.66%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.41'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.72'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("D31").Value = '0.0₃0725'
$ws.Range("E31").Value = '  -3.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.12'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.54'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.90%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("E35").Value = '  -0.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.26'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.02'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.908'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.71%  '
$ws.Range("E39").Value = '  -0.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '39.10'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '148.48'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.377'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.38%  '
$ws.Range("E43").Value = '  -1.92%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '280.84'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.15%  '
$ws.Range("E45").Value = '  -4.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0930'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0505'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.558'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.52%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.56'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.82%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0217'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.89%  '
$ws.Range("E51").Value = '  -1.43%  '
